# Update fig3 qp plot: insert a new "Sim-QAA" data column (column E) into
# each of the three worksheets, shifting the existing TNC/SNOPT/MATLAB/QCQP/
# IPOPT columns one place to the right (E->F, F->G, G->H, H->I, I->J).

$wb = $excel.ActiveWorkbook

# New column-E values (row -> value) for each worksheet, keyed by sheet index.
$values1 = @{2="0.222"; 3="0.249"; 4="0.267"; 5="0.535"; 6="0.346"; 7="0.32"; 8="0.313"; 9="0.342"; 10="0.413"; 11="0.288"}
$values2 = @{2="0.0003369736671447754"; 3="0.0002719881534576416"; 4="0.0003369717597961426"; 5="0.0003339846134185791"; 6="0.0003610746860504151"; 7="0.0003179309368133545"; 8="0.0003499839305877686"; 9="0.0003230266571044922"; 10="0.0003719723224639892"; 11="0.0003407449722290039"}
$values3 = @{2="0.006402499675750732"; 3="0.004623798608779907"; 4="0.005054576396942139"; 5="0.002337892293930054"; 6="0.003971821546554566"; 7="0.003815171241760254"; 8="0.004549791097640992"; 9="0.003876319885253906"; 10="0.003347750902175903"; 11="0.004770429611206055"}

$sheetValues = @($values1, $values2, $values3)

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Insert a new blank column before column E; this shifts old E:I -> F:J
    # and copies formatting/styles along with the existing data.
    $ws.Columns("E:E").Insert()

    # New header label for the inserted column.
    $ws.Range("E1").Value = "Sim-QAA"

    $rowValues = $sheetValues[$i - 1]
    foreach ($r in 2..11) {
        $ws.Cells.Item($r, 5).Value = [double]$rowValues["$r"]
    }
}

Write-Output "fig3 qp plot updated"
